$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Add($ws2.Range("I6"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/4cf339f51b8a2e9653f66b5437e6e7c9776585bb/e2e/22de129e-16a6-4966-b14d-39f5b6532b64.md", "", "", "22de129e-16a6-4966-b14d-39f5b6532b64.md")
$ws2.Range("I6").Style = "HyperLink"
